$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{
        A=@{T="n";V=112363550};
        B=@{T="n";V=5135};
        C=@{T="s";V="Ovaliderad"};
        D=@{T="s";V="LC"};
        E=@{T="n";V=105930};
        F=@{T="s";V="Vågbandad barkbock"};
        G=@{T="s";V="Semanotus undatus"};
        H=@{T="s";V="(Linnaeus, 1758)"};
        I=@{T="empty"};
        M=@{T="s";V="färska gnagspår"};
        P=@{T="s";V="Suderskog, Gtl"};
        Q=@{T="n";V=721928};
        R=@{T="n";V=6397835};
        S=@{T="n";V=10};
        T=@{T="s";V="Gotland"};
        U=@{T="s";V="Gotland"};
        V=@{T="s";V="Gotland"};
        W=@{T="s";V="Bäl"};
        Y=@{T="s";V="2023-09-22"};
        AA=@{T="s";V="2023-09-22"};
        AD=@{T="b";V=$false};
        AE=@{T="b";V=$false};
        AG=@{T="b";V=$false};
        AT=@{T="empty"};
        AW=@{T="s";V="Per Karlsson Linderum"};
        AX=@{T="s";V="Per Karlsson Linderum"};
        AY=@{T="empty"};
    };
    @{
        A=@{T="n";V=112363391};
        B=@{T="n";V=73758};
        C=@{T="s";V="Ovaliderad"};
        D=@{T="s";V="LC"};
        E=@{T="n";V=6426};
        F=@{T="s";V="Kattfotslav"};
        G=@{T="s";V="Felipes leucopellaeus"};
        H=@{T="s";V="(Ach.) Frisch & G.Thor"};
        I=@{T="empty"};
        J=@{T="s";V="bålar"};
        P=@{T="s";V="Suderskog, Gtl"};
        Q=@{T="n";V=721924};
        R=@{T="n";V=6397845};
        S=@{T="n";V=10};
        T=@{T="s";V="Gotland"};
        U=@{T="s";V="Gotland"};
        V=@{T="s";V="Gotland"};
        W=@{T="s";V="Bäl"};
        Y=@{T="s";V="2023-09-22"};
        AA=@{T="s";V="2023-09-22"};
        AD=@{T="b";V=$false};
        AE=@{T="b";V=$false};
        AG=@{T="b";V=$false};
        AT=@{T="empty"};
        AW=@{T="s";V="Per Karlsson Linderum"};
        AX=@{T="s";V="Per Karlsson Linderum"};
        AY=@{T="empty"};
    };
    @{
        A=@{T="n";V=112363523};
        B=@{T="n";V=85238};
        C=@{T="s";V="Ovaliderad"};
        D=@{T="s";V="LC"};
        E=@{T="n";V=3712};
        F=@{T="s";V="Blå slemspindling"};
        G=@{T="s";V="Cortinarius salor"};
        H=@{T="s";V="Fr."};
        I=@{T="empty"};
        J=@{T="s";V="fruktkroppar"};
        P=@{T="s";V="Suderskog, Gtl"};
        Q=@{T="n";V=721925};
        R=@{T="n";V=6397844};
        S=@{T="n";V=10};
        T=@{T="s";V="Gotland"};
        U=@{T="s";V="Gotland"};
        V=@{T="s";V="Gotland"};
        W=@{T="s";V="Bäl"};
        Y=@{T="s";V="2023-09-22"};
        AA=@{T="s";V="2023-09-22"};
        AD=@{T="b";V=$false};
        AE=@{T="b";V=$false};
        AG=@{T="b";V=$false};
        AT=@{T="empty"};
        AW=@{T="s";V="Per Karlsson Linderum"};
        AX=@{T="s";V="Per Karlsson Linderum"};
        AY=@{T="empty"};
    };
    @{
        A=@{T="n";V=112363369};
        B=@{T="n";V=93539};
        C=@{T="s";V="Ovaliderad"};
        D=@{T="s";V="LC"};
        E=@{T="n";V=2180};
        F=@{T="s";V="Blåmossa"};
        G=@{T="s";V="Leucobryum glaucum"};
        H=@{T="s";V="(Hedw.) Ångstr."};
        I=@{T="empty"};
        J=@{T="s";V="plantor/tuvor"};
        P=@{T="s";V="Suderskog, Gtl"};
        Q=@{T="n";V=721921};
        R=@{T="n";V=6397800};
        S=@{T="n";V=10};
        T=@{T="s";V="Gotland"};
        U=@{T="s";V="Gotland"};
        V=@{T="s";V="Gotland"};
        W=@{T="s";V="Bäl"};
        Y=@{T="s";V="2023-09-22"};
        AA=@{T="s";V="2023-09-22"};
        AD=@{T="b";V=$false};
        AE=@{T="b";V=$false};
        AG=@{T="b";V=$false};
        AT=@{T="empty"};
        AW=@{T="s";V="Per Karlsson Linderum"};
        AX=@{T="s";V="Per Karlsson Linderum"};
        AY=@{T="empty"};
    };
)
$startRow = 4
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $data = $rows[$i]
    foreach ($col in $data.Keys) {
        $cell = $ws.Range("$col$r")
        $info = $data[$col]
        switch ($info.T) {
            "n" {
                $cell.Value = $info.V
            }
            "b" {
                $cell.Value = $info.V
            }
            "s" {
                # Force text storage even for date-like strings (e.g. "2023-09-22")
                # so Excel doesn't auto-convert them to a date serial number.
                $cell.NumberFormat = "@"
                $cell.Value = $info.V
                $cell.ClearFormats()
            }
            "empty" {
                # Placeholder cell that exists but carries no value (matches the
                # empty inline-string cells in the source export). A harmless
                # no-op format touch is enough to make Excel materialize the
                # (empty) cell without leaving any real formatting behind.
                $cell.Font.Bold = $false
            }
        }
    }
}
